$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G1").Value = "Replication ES.r"
$ws.Range("D1").Value = "Median Replication ES.r"
$ws.Range("B1").Value = "ES.o"
$ws.Range("C1").Value = "95% CI Lower, Upper.o"
$ws.Range("H1").Value = "99% CI Lower, Upper.r"
$ws.Range("E1").Value = "Replication ES.r.unweighted.DNU"
$ws.Range("F1").Value = "99% CI Lower, Upper.r.unweighted.DNU"
$ws.Range("L1").Value = "Key statistics.r"
$ws.Range("M1").Value = "df.r"
$ws.Range("N1").Value = "N.r"
$ws.Range("O1").Value = "p.r"
$ws.Rows.AutoFit()
